{"js": "// Replace \"generally\" with \"not always\" in the sentence about teams that\n// resign their own players, i.e.\n//   \"... of the league) do generally see higher win rates ...\"\n// becomes\n//   \"... of the league) do not always see higher win rates ...\"\n\nconst body = context.document.body;\n\n// Search for the unique phrase spanning the edit point so we don't\n// accidentally touch any other occurrence of \"generally\" in the document.\nconst results = body.search(\" of the league) do generally see higher\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\" of the league) do not always see higher\", \"Replace\");\nawait context.sync();\n", "ps1": "# Replace \"generally\" with \"not always\" in the sentence about teams that\n# resign their own players, i.e.\n#   \"... of the league) do generally see higher win rates ...\"\n# becomes\n#   \"... of the league) do not always see higher win rates ...\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"of the league) do generally see higher\"\n$find.Replacement.Text = \"of the league) do not always see higher\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$result = $find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n)\n\nif (-not $result) {\n    throw \"Target phrase not found\"\n}\n"}
